# Apply the data-update edit to the "Import" worksheet.
#
# Rows 2-5 get replaced with new rows of data (new dates / values),
# and the previous row 6 values are cleared out (row shrinks from 5 data
# rows back down to 4).
#
# Finally the active selection is moved to F5, matching the new
# selection recorded in the sheet's sheetView.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-5 (Datum, Uhrzeit, Summe, Summe Aachen, Summe Todesfaelle, Summe genesen, Akute Faelle)
$data = @(
    @(44166, 0.46875,            10316, 4441, 184, 8968, 1164),
    @(44167, 0.42708333333333331, 10460, 4487, 188, 9194, 1078),
    @(44168, 0.41666666666666669, 10605, 4550, 191, 9314, 1100),
    @(44169, 0.41666666666666669, 10764, 4627, 196, 9548, 1020)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = 2 + $i
    $values = $data[$i]
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
    $ws.Cells.Item($row, 5).Value = $values[4]
    $ws.Cells.Item($row, 6).Value = $values[5]
    $ws.Cells.Item($row, 7).Value = $values[6]
}

# Clear out the old row 6 values (no new data available for that day yet)
$ws.Range("A6:G6").ClearContents()

# Update the active selection to match the new sheet view
$ws.Range("F5").Select()
